# Bill of Materials update.
# Adds a new BOM line item (IDC cable, row 11) and tidies up a couple of
# cosmetic leftovers (two cells that had stray formatting, page margins,
# and the header/footer text codes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 11: IDC cable part -------------------------------------------
$ws.Range("A11").Value = 5
$ws.Range("B11").Value = "H3CCH-1018G"
$ws.Range("C11").Value = "Assmann WSW Components"
$ws.Range("F11").Value = "IDC CBL - HHKC10H/AE10G/HHKC10H"

# --- Clear the stray formatting that F3/F5 carried -------------------------
$ws.Range("F3").Font.Name = "Arial"
$ws.Range("F3").Font.Size = 10
$ws.Range("F5").Font.Name = "Arial"
$ws.Range("F5").Font.Size = 10

# --- Page margins (top/bottom tightened slightly) ---------------------------
$ws.PageSetup.TopMargin = 73.8
$ws.PageSetup.BottomMargin = 73.8

# --- Header/footer: drop the explicit Times New Roman 12pt codes -----------
$ws.PageSetup.CenterHeader = "&A"
$ws.PageSetup.CenterFooter = "Page &P"

# --- Move the active selection to the new last row --------------------------
$ws.Range("B11").Select() | Out-Null
